$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Insert a new data row at row 12 for "OPLEX-N SYRUP 125ML"
#    (rows 12-16 shift down to 13-17)
# ------------------------------------------------------------------
$ws.Rows.Item(12).Insert()

# Re-apply the standard item-row formatting (styles s=7..12 over A:Q)
# by copying the format from the row directly above (row 11, which is
# still a normal item row) onto the freshly inserted blank row 12.
$ws.Range("A11:Q11").Copy()
$ws.Range("A12:Q12").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Match the row height used by the other item rows (alternates
# 25.5 / 24.75 - row 12 takes the same height the old row 12 had).
$ws.Rows.Item(12).RowHeight = 25.5

# Re-create the merges for row 12 that mirror every other item row.
$ws.Range("A12:B12").Merge()
$ws.Range("C12:G12").Merge()
$ws.Range("H12:K12").Merge()
$ws.Range("L12:M12").Merge()
$ws.Range("N12:O12").Merge()

# ------------------------------------------------------------------
# 2. Fill in the values for the new row
# ------------------------------------------------------------------
$ws.Range("A12").Value = 6
$ws.Range("C12").Value = "OPLEX-N SYRUP 125ML"
$ws.Range("H12").Value = "4:0"
$ws.Range("L12").Value = "'1"
$ws.Range("N12").Value = "31.00"
$ws.Range("P12").Value = "'31.0000"
$ws.Range("Q12").Value = "1:0"

# Fix up the text-that-looks-numeric cells (L12, P12) so they keep the
# *same* cell style as their neighbours instead of the ad-hoc
# "quoted-text" style Excel assigns when a value is typed with a
# leading apostrophe.
$ws.Range("L11").Copy()
$ws.Range("L12").PasteSpecial(-4122)
$ws.Range("P11").Copy()
$ws.Range("P12").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ------------------------------------------------------------------
# 3. Renumber the "م" (index) column for the rows that followed the
#    inserted row - they have each shifted down by one position.
# ------------------------------------------------------------------
$ws.Range("A13").Value = 7
$ws.Range("A14").Value = 8
$ws.Range("A15").Value = 9

# ------------------------------------------------------------------
# 4. Update the totals row (now row 16) with the new grand total.
# ------------------------------------------------------------------
$ws.Range("P16").Value = 513.3
$ws.Rows.Item(16).RowHeight = 25.5

# ------------------------------------------------------------------
# 5. Update the footer timestamp (now row 17, column A).
# ------------------------------------------------------------------
$ws.Range("A17").Value = "Sunday, 27 July, 2025 10:34 AM"

Write-Output "edit applied"
